$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph (the 2nd paragraph, right after
#    the Heading1 title). It contained a bold "Meta description" run plus a
#    plain run with the description text - the whole paragraph (and its
#    paragraph mark) goes away.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 2. Insert a new paragraph just before the final ("Prompt: ...") paragraph
#    that carries the bold title text which used to live in the meta
#    description paragraph we just removed.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$precedingPara = $d.Paragraphs($n - 1)
$insertionPoint = $precedingPara.Range.Duplicate
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs($n)
$newRange = $newPara.Range.Duplicate
$newRange.Collapse(1)
$newRange.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Aloha Fruit Bonanza Free - Review &amp; Guide | RTP 97.01%</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# ---------------------------------------------------------------------------
# 3. Replace the "Prompt: ..." text with the new description text, keeping
#    the italic formatting already present on that run.
# ---------------------------------------------------------------------------
$oldPrompt = "Prompt: Create a cartoon-style feature image for the game " + [char]34 + "Aloha Fruit Bonanza" + [char]34 + ". The image should feature a happy Maya warrior wearing glasses. The design should have a tropical feel, with bright colors and a beach background. It should also prominently feature fruit symbols from the game, such as watermelon, coconut cocktail, and dragon fruit cocktail, as well as the red number seven. The warrior should be holding a slot machine lever, and there should be cascading symbols falling around him. The overall feeling of the image should be fun and exciting, capturing the lightheartedness of the game."
$newPrompt = "Discover the stunning visuals and innovative Scatter Pays feature in Aloha Fruit Bonanza. With an RTP of 97.01%, play for free and win up to 7,500x your stake."

$d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newPrompt, 2) | Out-Null
